$d = $word.ActiveDocument

# Locate the target paragraph: "... 'Neue datei einstellungen' moeglicherweise
# nicht nur temporaer sondern auch lokal oder so speichern". It is the final
# paragraph in the document (right before the sectPr).
$target = $d.Paragraphs.Last
if (-not ($target.Range.Text -like "*sondern auch lokal*")) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*sondern auch lokal*") {
            $target = $p
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

# Replace the paragraph with:
#  1) the same paragraph, but with its last run split in three so that
#     "temporär" is wrapped in a gramStart/gramEnd proofing-error pair
#     (mirrors what Word's grammar checker inserts when the text is re-typed), and
#  2) a new list-item paragraph right after it: "Convert app.py to a class…"
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="4FCB7B8B" w14:textId="5A9442FA" w:rsidR="00B53FAB" w:rsidRPr="007A7336" w:rsidRDefault="002E7383" w:rsidP="00305483"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">„Neue </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>datei</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>einstellungen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">“ möglicherweise nicht nur </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>temporär</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> sondern auch lokal oder so speichern</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Convert app.py to a c</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>lass…</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$target.Range.InsertXML($xml) | Out-Null
